$wb = $excel.ActiveWorkbook

# --- Extensions(Ctrl+Shift+X) sheet: add a new "Azure App Service" extension row ---
$wsExt = $wb.Worksheets.Item("Extensions(Ctrl+Shift+X)")
[void]$wsExt.Activate()
$newRow = $wsExt.Cells.Item(15, 1)
$newRow.Value = "Azure App Service" + [char]10 + "ms-azuretools.vscode-azureappservice"
$newRow.WrapText = $true
[void]$wsExt.Range("A16").Select()

# --- DatingApp sheet: add PowerShell execution-policy command next to the dotnet publish command ---
$wsApp = $wb.Worksheets.Item("DatingApp")
[void]$wsApp.Activate()
$wsApp.Range("L2").Value = "Set-ExecutionPolicy RemoteSigned -Scope CurrentUser"
[void]$wsApp.Range("L3").Select()
